$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HomeResource")

# Fix the "Carma" (Karma) row's name text stays the same; no content change needed there
# (it was only a shared-string index shuffle caused by the other edits)

# Row 3 (food/Hearth): add detail text in column F
$ws.Range("F3").Value = "Lugar para armazenar suprimentos"

# Row 4 (money/Gold): add detail text in column F
$ws.Range("F4").Value = "Espaço onde os moradores vivem"

# Row 5 (knowledge): replace "(Indisponível)" with "(Não implementado)"
$ws.Range("F5").Value = "(Não implementado)"

# Row 6 (culture): replace "(Indisponível)" with "(Não implementado)"
$ws.Range("F6").Value = "(Não implementado)"

# Row 7 (faith): add detail text in column F
$ws.Range("F7").Value = "Lugar para cultivar e colher plantações"

# Row 8 (fun): replace "(Indisponível)" with "(Não implementado)"
$ws.Range("F8").Value = "(Não implementado)"

# Row 9 (medicine): replace "(Indisponível)" with "(Não implementado)"
$ws.Range("F9").Value = "(Não implementado)"

# Update column F width to fit new (wider) content
$ws.Columns.Item(6).ColumnWidth = 35.166666666666664

# Move selection to F9, matching the final saved cursor position
$ws.Range("F9").Select()
